$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") entirely; D and E shift left to C and D.
$ws.Columns.Item(3).Delete()

# Update header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Update data row
$ws.Range("B2").Value = -27.19826829767476
$ws.Range("C2").Value = "s__Amedibacillus dolichus"
$ws.Range("D2").Value = "s__Amedibacillus dolichus(reject)"
